# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "HUBER ARLEY GIL GOMEZ" / periodo 2505 row (row 18).
# Everything below it (the 2504 row, the signature block, etc.) shifts up by one row.
$ws.Range("B18:J18").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

# Update the summary values at the top of the sheet.
$ws.Range("E11").Value = 341642
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 3

# Rows 16 and 17 (LUIS ALBERTO AREVALO ARENILLA) swap which period/amount they show.
$ws.Range("E16").Value = "2506"
$ws.Range("F16").Value = 85410
$ws.Range("E17").Value = "2507"
$ws.Range("F17").Value = 128116

# Row 18 now holds the data that used to be the old "2504" row (after the shift-up above).
# Replace it with the new worker entry: LUIS ALBERTO AREVALO ARENILLA, periodo 2508.
$ws.Range("C18").Value = "19772276"
$ws.Range("D18").Value = "LUIS ALBERTO AREVALO ARENILLA"
$ws.Range("E18").Value = "2508"
$ws.Range("F18").Value = 128116
$ws.Range("G18").Value = 3202875

# Persist the changes so the workbook recomputes/compacts its internal string table.
$wb.Save()
